# Added code for Multi Doc Upload functionality
$wb = $excel.ActiveWorkbook

# --- Update selection on the existing "Templates" sheet -------------------
# Select the whole header row (matches the new workbook selection state,
# and causes "tabSelected" to move off this sheet once another sheet
# becomes active below).
$wsTemplates = $wb.Worksheets.Item("Templates")
$wsTemplates.Rows.Item(1).Select() | Out-Null

# --- Add the new "MultiDocs" sheet after "Templates" -----------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMultiDocs = $wb.Worksheets.Add($null, $lastSheet)
$wsMultiDocs.Name = "MultiDocs"

# Header row
$wsMultiDocs.Cells.Item(1,1).Value = "folder"
$wsMultiDocs.Cells.Item(1,2).Value = "contact"
$wsMultiDocs.Cells.Item(1,3).Value = "company"
$wsMultiDocs.Cells.Item(1,4).Value = "prospect"
$wsMultiDocs.Cells.Item(1,5).Value = "task"
$wsMultiDocs.Cells.Item(1,6).Value = "case"
$wsMultiDocs.Cells.Item(1,7).Value = "tags"

# Data rows, populated column by column
$wsMultiDocs.Cells.Item(2,1).Value = "Root directory"
$wsMultiDocs.Cells.Item(3,1).Value = "Word Templates"

$wsMultiDocs.Cells.Item(2,2).Value = "contact1"
$wsMultiDocs.Cells.Item(3,2).Value = "contact2"

$wsMultiDocs.Cells.Item(2,3).Value = "company1"
$wsMultiDocs.Cells.Item(3,3).Value = "company2"

$wsMultiDocs.Cells.Item(2,4).Value = "prospect1"
$wsMultiDocs.Cells.Item(3,4).Value = "prospect2"

$wsMultiDocs.Cells.Item(2,5).Value = "task1"
$wsMultiDocs.Cells.Item(3,5).Value = "task2"

$wsMultiDocs.Cells.Item(2,6).Value = "case1"
$wsMultiDocs.Cells.Item(3,6).Value = "case2"

$wsMultiDocs.Cells.Item(2,7).Value = "tag1"
$wsMultiDocs.Cells.Item(3,7).Value = "tag2"

# Header fill (same yellow highlight style used by the other sheets)
$wsMultiDocs.Range("A1:G1").Interior.Color = 65535

# Column A width, best-fit like the other sheets
$wsMultiDocs.Columns.Item(1).ColumnWidth = 13.73

# Make MultiDocs the active sheet/selection, matching the target view state
$wsMultiDocs.Range("H3").Select() | Out-Null
